$wb = $excel.ActiveWorkbook

function Set-TextValue($ws, $cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

$wsFeatures = $wb.Worksheets.Item("Features")
$wsGlobal = $wb.Worksheets.Item("Global Metrics")

Set-TextValue $wsFeatures "B2" "0,583"
Set-TextValue $wsFeatures "C2" "0,318"
Set-TextValue $wsFeatures "D2" "0,412"
Set-TextValue $wsFeatures "E2" "0,969"
Set-TextValue $wsFeatures "B3" "0,667"
Set-TextValue $wsFeatures "C3" "0,364"
Set-TextValue $wsFeatures "D3" "0,471"
Set-TextValue $wsFeatures "E3" "0,828"
Set-TextValue $wsFeatures "B4" "0,417"
Set-TextValue $wsFeatures "C4" "0,227"
Set-TextValue $wsFeatures "D4" "0,294"
Set-TextValue $wsFeatures "E4" "0,956"
Set-TextValue $wsFeatures "B5" "0,500"
Set-TextValue $wsFeatures "C5" "0,286"
Set-TextValue $wsFeatures "D5" "0,364"
Set-TextValue $wsFeatures "E5" "0,915"
Set-TextValue $wsFeatures "B6" "0,375"
Set-TextValue $wsFeatures "C6" "0,150"
Set-TextValue $wsFeatures "D6" "0,214"
Set-TextValue $wsFeatures "E6" "0,895"
Set-TextValue $wsFeatures "B7" "0,250"
Set-TextValue $wsFeatures "C7" "0,100"
Set-TextValue $wsFeatures "D7" "0,143"
Set-TextValue $wsFeatures "E7" "0,937"
Set-TextValue $wsFeatures "E8" "0,691"
Set-TextValue $wsFeatures "B9" "0,200"
Set-TextValue $wsFeatures "C9" "0,083"
Set-TextValue $wsFeatures "D9" "0,118"
Set-TextValue $wsFeatures "E9" "0,696"
Set-TextValue $wsFeatures "B10" "0,667"
Set-TextValue $wsFeatures "C10" "0,182"
Set-TextValue $wsFeatures "D10" "0,286"
Set-TextValue $wsFeatures "E10" "0,797"
Set-TextValue $wsFeatures "B11" "0,600"
Set-TextValue $wsFeatures "C11" "0,462"
Set-TextValue $wsFeatures "D11" "0,522"
Set-TextValue $wsFeatures "E11" "0,947"
Set-TextValue $wsFeatures "B12" "0,500"
Set-TextValue $wsFeatures "C12" "0,385"
Set-TextValue $wsFeatures "D12" "0,435"
Set-TextValue $wsFeatures "E12" "0,950"
Set-TextValue $wsFeatures "B13" "0,400"
Set-TextValue $wsFeatures "C13" "0,333"
Set-TextValue $wsFeatures "D13" "0,364"
Set-TextValue $wsFeatures "E13" "0,901"
Set-TextValue $wsFeatures "B14" "0,300"
Set-TextValue $wsFeatures "C14" "0,250"
Set-TextValue $wsFeatures "D14" "0,273"
Set-TextValue $wsFeatures "E14" "0,945"
Set-TextValue $wsFeatures "B15" "1,000"
Set-TextValue $wsFeatures "C15" "0,231"
Set-TextValue $wsFeatures "D15" "0,375"
Set-TextValue $wsFeatures "E15" "0,231"
Set-TextValue $wsFeatures "B16" "0,200"
Set-TextValue $wsFeatures "C16" "0,100"
Set-TextValue $wsFeatures "D16" "0,133"
Set-TextValue $wsFeatures "E16" "0,892"
Set-TextValue $wsFeatures "B17" "0,429"
Set-TextValue $wsFeatures "C17" "0,273"
Set-TextValue $wsFeatures "D17" "0,333"
Set-TextValue $wsFeatures "E17" "0,872"
Set-TextValue $wsFeatures "E18" "0,519"
Set-TextValue $wsFeatures "B19" "0,667"
Set-TextValue $wsFeatures "C19" "0,200"
Set-TextValue $wsFeatures "D19" "0,308"
Set-TextValue $wsFeatures "E19" "0,391"
Set-TextValue $wsFeatures "B20" "0,400"
Set-TextValue $wsFeatures "C20" "0,250"
Set-TextValue $wsFeatures "D20" "0,308"
Set-TextValue $wsFeatures "E20" "0,784"
Set-TextValue $wsFeatures "E21" "0,699"
Set-TextValue $wsFeatures "B28" "1,000"
Set-TextValue $wsFeatures "C28" "1,000"
Set-TextValue $wsFeatures "D28" "1,000"
Set-TextValue $wsFeatures "E28" "1,000"
Set-TextValue $wsFeatures "E29" "0,964"
Set-TextValue $wsFeatures "B32" "1,000"
Set-TextValue $wsFeatures "C32" "0,400"
Set-TextValue $wsFeatures "D32" "0,571"
Set-TextValue $wsFeatures "E32" "0,400"
Set-TextValue $wsFeatures "E38" "1,000"
Set-TextValue $wsGlobal "B2" "0,128"
Set-TextValue $wsGlobal "C2" "0,355"
Set-TextValue $wsGlobal "D2" "0,236"
Set-TextValue $wsGlobal "E2" "0,651"
